# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.463.06"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.672.78"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'617.99"
$ws.Range("E5").Value = "  -8.20%  "
$ws.Range("D6").Value = "'159.32"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D8").Value = "'0.497"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").Value = "4.290.28"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'32.43"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "3.671.38"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "69.496.32"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "'6.51"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'15.89"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("D21").Value = "'469.05"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D23").Value = "'79.42"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "3.817.61"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'11.10"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  -4.56%  "
$ws.Range("D28").Value = "'8.67"
$ws.Range("E28").Value = "  -5.40%  "
$ws.Range("D29").Value = "'2.61"
$ws.Range("E29").Value = "  -3.34%  "
$ws.Range("D30").Value = "'1.67"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").Value = "'26.62"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").Value = "'6.39"
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("D36").Value = "3.671.06"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "'8.29"
$ws.Range("D39").Value = "'178.68"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.23"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'5.80"
$ws.Range("E42").Value = "  -5.35%  "
$ws.Range("D43").Value = "'0.0893"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").Value = "'0.926"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'29.19"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'46.78"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "'0.000265"
$ws.Range("E49").Value = "  -5.97%  "
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("E51").Value = "  -7.01%  "
